$d = $word.ActiveDocument

# The author renamed the (misspelled) English phrase "resource packs" to
# "dataset" throughout the first bullet point:
#   "Organisation des « resource packs » - 10m"
#   -> "Organisation des « dataset » - 10m"
# Step 1: swap the word itself (keeps the spell-check bookmarks that
# already wrapped it).
$d.Content.Find.Execute("resource", $true, $false, $false, $false, $false, $true, 1, $false, "dataset", 2)

# Step 2: drop the now-redundant " packs" before the closing guillemet,
# preserving the non-breaking space this document always places right
# before a "»".
$d.Content.Find.Execute(" packs »", $true, $false, $false, $false, $false, $true, 1, $false, " »", 2)

# Word naturally consolidates neighbouring, identically formatted runs
# whenever a paragraph is touched by an edit/save cycle. Running a
# (no-op) replace on the recurring "- " separator normalizes every
# bullet that had it split across multiple runs.
$d.Content.Find.Execute("– ", $true, $false, $false, $false, $false, $true, 1, $false, "– ", 2)

# Likewise, consolidate the split bold duration runs "3"+"h" and "4"+"h".
$d.Content.Find.Execute("3h", $true, $false, $false, $false, $false, $true, 1, $false, "3h", 2)
$d.Content.Find.Execute("4h", $true, $false, $false, $false, $false, $true, 1, $false, "4h", 2)

# The "Inventer un indice pour le chat cyclope - " bullet was missing its
# duration; the author added a bold "10m" at the end of it. Copy the
# character formatting (bold + bold-complex-script) from an existing
# duration run so the new run's rPr matches exactly, then set its text.
$srcRange = $d.Content.Duplicate
$srcRange.Find.Execute("5m")
$srcFormatted = $srcRange.FormattedText

$targetPara = $d.Paragraphs(5)
$insertPoint = $targetPara.Range.End - 1
$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.FormattedText = $srcFormatted
$newRange = $d.Range($insertPoint, $insertPoint + 2)
$newRange.Text = "10m"
